# Add macro support: a new "macros" worksheet holding token/expansion
# pairs, plus updates to the Artist column on "Final" so a few cards
# reference those macro tokens via {{ }} placeholders.

$wb = $excel.ActiveWorkbook

# --- Add the new "macros" worksheet, placed right after "Final" ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "macros"

$final = $wb.Worksheets.Item("Final")
$newSheet.Move($null, $final)

# Re-fetch a live handle to the freshly positioned sheet before writing.
$macros = $wb.Worksheets.Item("macros")
$macros.Range("A1").Value = "FOOP"
$macros.Range("A2").Value = "DOOP"
$macros.Range("A3").Value = "BLOOP"
$macros.Range("B1").Value = "foop"
$macros.Range("B2").Value = "<`$b>doop<`$>"
$macros.Range("B3").Value = "bloop bloop bloop"

# --- Point a few cards' Artist values at the new macro tokens ---
$ws = $wb.Worksheets.Item("Final")
$ws.Range("O13").Value = "Arty {{ DOOP }}"
$ws.Range("O12").Value = "{{ FOOP }} McArtison"
$ws.Range("O11").Value = "Arty {{ BLOOP }} McArtison"

# Leave the view roughly where the author left it.
$ws.Activate()
$ws.Range("P20").Select()
